$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BLSD_Students")

# Force the whole row to stay stored as text, matching the original
# shared-string (t="s") cell types for A2:H2 -- otherwise Excel will
# auto-coerce the date/number-looking strings into real date/number cells.
$ws.Range("A2:H2").NumberFormat = "@"

$ws.Range("A2").Value = "cbd81bf4-58cd-4671-931d-9e712a51bb3e"
$ws.Range("B2").Value = "Q12312312Q"
$ws.Range("C2").Value = "Eda"
$ws.Range("D2").Value = "Isaku"
$ws.Range("E2").Value = "eda@gmail.com"
$ws.Range("F2").Value = "BLSD"
$ws.Range("G2").Value = "2022-04-22"
$ws.Range("H2").Value = "3"
